# Add a new vocabulary entry ("den Schlüssel hinterlegen / abgeben" / "to leave
# the key") to the "vocab" sheet, inserted right after the existing row for
# "anwesend / abwesend" / "present / absent" (row 268), pushing the rows below
# it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vocab")
$ws.Activate()

# Insert a new row at 269; this shifts rows 269:292 down to 270:293 and
# copies the number formatting (short date in column C) from the
# neighbouring rows automatically.
$ws.Rows("269").Insert()

# German term (column A). The original author's text has mixed run
# formatting (an explicit Calibri/size-11 run starting after "den Schl"),
# so re-create that via Characters() once the plain value is in place.
$ws.Range("A269").Value = "den Schlüssel hinterlegen / abgeben"
$runChars = $ws.Range("A269").Characters(9, 27)
$runChars.Font.Name = "Calibri"
$runChars.Font.Size = 11

# English translation (column B).
$ws.Range("B269").Value = "to leave the key"

# Lesson date (column C) - same lesson date as the surrounding rows.
$ws.Range("C269").Value = 44657

# Lesson number (column D) - same lesson number as the surrounding rows.
$ws.Range("D269").Value = 10

# Phrase / Word category (column E) - this entry is tagged as a "word".
$ws.Range("E269").Value = "word"

# Update the on-screen selection to match the author's saved view.
$ws.Range("A274").Select()
